$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H33").Value = 1184.5333
$ws_ALC.Range("I33").Value = 1229.0769
$ws_ALC.Range("K33").Value = 1229.0769
$ws_ALC.Range("M33").Value = -1000.0769
$ws_ALC.Range("H38").Value = 69
$ws_ALC.Range("I38").Value = 69
$ws_ALC.Range("K38").Value = 207
$ws_ALC.Range("M38").Value = 165
$ws_ALC.Range("H40").Value = 3708.0908
$ws_ALC.Range("I40").Value = 3483.3333
$ws_ALC.Range("J40").Value = 3977.8
$ws_ALC.Range("K40").Value = 3483.3333
$ws_ALC.Range("L40").Value = 3977.8
$ws_ALC.Range("M40").Value = -3308.3333
$ws_ALC.Range("N40").Value = -4327.8
$ws_ALC.Range("H41").Value = 9618170
$ws_ALC.Range("I41").Value = 13892479
$ws_ALC.Range("J41").Value = 974.75
$ws_ALC.Range("K41").Value = 13892479
$ws_ALC.Range("L41").Value = 974.75
$ws_ALC.Range("M41").Value = -13892039
$ws_ALC.Range("N41").Value = -1854.75
$ws_ALC.Range("H132").Value = 2724.1052
$ws_ALC.Range("I132").Value = 2839.9285
$ws_ALC.Range("K132").Value = 8519.7855
$ws_ALC.Range("M132").Value = -5989.7855
$ws_ALC.Range("H137").Value = 2290.0356
$ws_ALC.Range("I137").Value = 2246.4
$ws_ALC.Range("K137").Value = 6739.200000000001
$ws_ALC.Range("M137").Value = -4189.200000000001
$ws_ALC.Range("H138").Value = 4561.1714
$ws_ALC.Range("I138").Value = 1690.375
$ws_ALC.Range("J138").Value = 5411.778
$ws_ALC.Range("K138").Value = 5071.125
$ws_ALC.Range("L138").Value = 16235.334
$ws_ALC.Range("M138").Value = 68.875
$ws_ALC.Range("N138").Value = -26515.334
$ws_ALC.Range("H141").Value = 3996.25
$ws_ALC.Range("I141").Value = 3995.2
$ws_ALC.Range("J141").Value = 3998
$ws_ALC.Range("K141").Value = 11985.6
$ws_ALC.Range("L141").Value = 11994
$ws_ALC.Range("M141").Value = -6805.599999999999
$ws_ALC.Range("N141").Value = -22354

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 2254.5
$ws_ARM.Range("I32").Value = 2304.7605
$ws_ARM.Range("K32").Value = 2304.7605
$ws_ARM.Range("M32").Value = -2017.7605
$ws_ARM.Range("H61").Value = 5134.877
$ws_ARM.Range("I61").Value = 2786.587
$ws_ARM.Range("J61").Value = 14955
$ws_ARM.Range("K61").Value = 2786.587
$ws_ARM.Range("L61").Value = 14955
$ws_ARM.Range("M61").Value = -2574.587
$ws_ARM.Range("N61").Value = -15379
$ws_ARM.Range("H74").Value = 70604.22
$ws_ARM.Range("I74").Value = 201108
$ws_ARM.Range("J74").Value = 5352.3335
$ws_ARM.Range("K74").Value = 201108
$ws_ARM.Range("L74").Value = 5352.3335
$ws_ARM.Range("M74").Value = -200234
$ws_ARM.Range("N74").Value = -7100.3335
$ws_ARM.Range("H77").Value = 70604.22
$ws_ARM.Range("I77").Value = 201108
$ws_ARM.Range("J77").Value = 5352.3335
$ws_ARM.Range("K77").Value = 1005540
$ws_ARM.Range("L77").Value = 26761.6675
$ws_ARM.Range("M77").Value = -1001172
$ws_ARM.Range("N77").Value = -35497.6675
$ws_ARM.Range("H122").Value = 10185.483
$ws_ARM.Range("I122").Value = 12105.2
$ws_ARM.Range("K122").Value = 36315.60000000001
$ws_ARM.Range("M122").Value = -33865.60000000001
$ws_ARM.Range("H125").Value = 40771.133
$ws_ARM.Range("J125").Value = 40771.133
$ws_ARM.Range("L125").Value = 40771.133
$ws_ARM.Range("N125").Value = -50611.133
$ws_ARM.Range("H136").Value = 5134.877
$ws_ARM.Range("I136").Value = 2786.587
$ws_ARM.Range("J136").Value = 14955
$ws_ARM.Range("K136").Value = 8359.761
$ws_ARM.Range("L136").Value = 44865
$ws_ARM.Range("M136").Value = -5809.761
$ws_ARM.Range("N136").Value = -49965

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H76").Value = 23999.5
$ws_BSM.Range("J76").Value = 23999.5
$ws_BSM.Range("L76").Value = 23999.5
$ws_BSM.Range("N76").Value = -24629.5
$ws_BSM.Range("H79").Value = 23999.5
$ws_BSM.Range("J79").Value = 23999.5
$ws_BSM.Range("L79").Value = 23999.5
$ws_BSM.Range("N79").Value = -26183.5
$ws_BSM.Range("H99").Value = 7578311.5
$ws_BSM.Range("I99").Value = 2295.4443
$ws_BSM.Range("K99").Value = 2295.4443
$ws_BSM.Range("M99").Value = -797.4443000000001
$ws_BSM.Range("H107").Value = 56253190
$ws_BSM.Range("I107").Value = 86539330
$ws_BSM.Range("J107").Value = 7512.143
$ws_BSM.Range("K107").Value = 86539330
$ws_BSM.Range("L107").Value = 7512.143
$ws_BSM.Range("M107").Value = -86537410
$ws_BSM.Range("N107").Value = -11352.143
$ws_BSM.Range("H128").Value = 3416.125
$ws_BSM.Range("I128").Value = 3416.125
$ws_BSM.Range("K128").Value = 10248.375
$ws_BSM.Range("M128").Value = -7758.375

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 6381.22
$ws_CRP.Range("I31").Value = 2615.6287
$ws_CRP.Range("K31").Value = 2615.6287
$ws_CRP.Range("M31").Value = -2320.6287
$ws_CRP.Range("H34").Value = 6381.22
$ws_CRP.Range("I34").Value = 2615.6287
$ws_CRP.Range("K34").Value = 2615.6287
$ws_CRP.Range("M34").Value = -2413.6287
$ws_CRP.Range("H58").Value = 9095776
$ws_CRP.Range("I58").Value = 14287299
$ws_CRP.Range("J58").Value = 10611.8
$ws_CRP.Range("K58").Value = 14287299
$ws_CRP.Range("L58").Value = 10611.8
$ws_CRP.Range("M58").Value = -14287096
$ws_CRP.Range("N58").Value = -11017.8
$ws_CRP.Range("H105").Value = 3404611.8
$ws_CRP.Range("I105").Value = 4202726.5
$ws_CRP.Range("K105").Value = 4202726.5
$ws_CRP.Range("M105").Value = -4200979.5
$ws_CRP.Range("H132").Value = 5506.25
$ws_CRP.Range("I132").Value = 2566
$ws_CRP.Range("K132").Value = 7698
$ws_CRP.Range("M132").Value = -5168
$ws_CRP.Range("H134").Value = 4715.2417
$ws_CRP.Range("I134").Value = 2035.1316
$ws_CRP.Range("J134").Value = 8958.75
$ws_CRP.Range("K134").Value = 6105.3948
$ws_CRP.Range("L134").Value = 26876.25
$ws_CRP.Range("M134").Value = -3570.3948
$ws_CRP.Range("N134").Value = -31946.25
$ws_CRP.Range("H136").Value = 9095776
$ws_CRP.Range("I136").Value = 14287299
$ws_CRP.Range("J136").Value = 10611.8
$ws_CRP.Range("K136").Value = 42861897
$ws_CRP.Range("L136").Value = 31835.4
$ws_CRP.Range("M136").Value = -42859347
$ws_CRP.Range("N136").Value = -36935.39999999999
$ws_CRP.Range("H141").Value = 79124.69
$ws_CRP.Range("J141").Value = 79124.69
$ws_CRP.Range("L141").Value = 79124.69
$ws_CRP.Range("N141").Value = -89484.69

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H82").Value = 32666
$ws_CUL.Range("I82").Value = 18999.5
$ws_CUL.Range("K82").Value = 56998.5
$ws_CUL.Range("M82").Value = -56592.5
$ws_CUL.Range("H85").Value = 32666
$ws_CUL.Range("I85").Value = 18999.5
$ws_CUL.Range("K85").Value = 56998.5
$ws_CUL.Range("M85").Value = -55594.5
$ws_CUL.Range("H136").Value = 1010.8125
$ws_CUL.Range("I136").Value = 1010.8125
$ws_CUL.Range("K136").Value = 3032.4375
$ws_CUL.Range("M136").Value = 2067.5625
$ws_CUL.Range("H138").Value = 4687.2144
$ws_CUL.Range("I138").Value = 4172.3335
$ws_CUL.Range("J138").Value = 7776.5
$ws_CUL.Range("K138").Value = 12517.0005
$ws_CUL.Range("L138").Value = 23329.5
$ws_CUL.Range("M138").Value = -7377.000499999998
$ws_CUL.Range("N138").Value = -33609.5

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H70").Value = 10444.546
$ws_GSM.Range("I70").Value = 9717
$ws_GSM.Range("J70").Value = 11317.6
$ws_GSM.Range("K70").Value = 9717
$ws_GSM.Range("L70").Value = 11317.6
$ws_GSM.Range("M70").Value = -9447
$ws_GSM.Range("N70").Value = -11857.6
$ws_GSM.Range("H73").Value = 10444.546
$ws_GSM.Range("I73").Value = 9717
$ws_GSM.Range("J73").Value = 11317.6
$ws_GSM.Range("K73").Value = 9717
$ws_GSM.Range("L73").Value = 11317.6
$ws_GSM.Range("M73").Value = -8781
$ws_GSM.Range("N73").Value = -13189.6
$ws_GSM.Range("H122").Value = 2654567.8
$ws_GSM.Range("J122").Value = 2555.25
$ws_GSM.Range("L122").Value = 7665.75
$ws_GSM.Range("N122").Value = -12565.75
$ws_GSM.Range("H126").Value = 9999
$ws_GSM.Range("I126").Value = 0
$ws_GSM.Range("K126").Value = 0
$ws_GSM.Range("M126").ClearContents()
$ws_GSM.Range("H140").Value = 76694.75
$ws_GSM.Range("J140").Value = 76694.75
$ws_GSM.Range("L140").Value = 76694.75
$ws_GSM.Range("N140").Value = -87054.75
$ws_GSM.Range("H141").Value = 28247.75
$ws_GSM.Range("J141").Value = 32665.555
$ws_GSM.Range("L141").Value = 32665.555
$ws_GSM.Range("N141").Value = -43025.555

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H40").Value = 6090.5625
$ws_LTW.Range("I40").Value = 2888.3333
$ws_LTW.Range("J40").Value = 8011.9
$ws_LTW.Range("K40").Value = 2888.3333
$ws_LTW.Range("L40").Value = 8011.9
$ws_LTW.Range("M40").Value = -2752.3333
$ws_LTW.Range("N40").Value = -8283.9
$ws_LTW.Range("H46").Value = 4832946.5
$ws_LTW.Range("I46").Value = 761.0833
$ws_LTW.Range("K46").Value = 761.0833
$ws_LTW.Range("M46").Value = -573.0833
$ws_LTW.Range("H127").Value = 59112.25
$ws_LTW.Range("J127").Value = 59112.25
$ws_LTW.Range("L127").Value = 59112.25
$ws_LTW.Range("N127").Value = -69032.25

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H2").Value = 49332.668
$ws_WVR.Range("I2").Value = 49332.668
$ws_WVR.Range("K2").Value = 49332.668
$ws_WVR.Range("M2").Value = -49220.668
$ws_WVR.Range("H4").Value = 15000
$ws_WVR.Range("I4").Value = 0
$ws_WVR.Range("K4").Value = 0
$ws_WVR.Range("M4").ClearContents()
$ws_WVR.Range("H132").Value = 11120807
$ws_WVR.Range("I132").Value = 13892959
$ws_WVR.Range("J132").Value = 32200.223
$ws_WVR.Range("K132").Value = 41678877
$ws_WVR.Range("L132").Value = 96600.66900000001
$ws_WVR.Range("M132").Value = -41676347
$ws_WVR.Range("N132").Value = -101660.669
$ws_WVR.Range("H136").Value = 40006180
$ws_WVR.Range("I136").Value = 125001100
$ws_WVR.Range("K136").Value = 375003300
$ws_WVR.Range("M136").Value = -375000750
$ws_WVR.Range("H141").Value = 0
$ws_WVR.Range("J141").Value = 0
$ws_WVR.Range("L141").Value = 0
$ws_WVR.Range("N141").ClearContents()
